# Added tests for a spot being full
# Mark several spots as "full": bump the challenged-amount counter (column C)
# and, when a spot becomes completely full, flip its Open flag (column B) to TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spot "B6" (row 21): one more challenger, still not full
$ws.Range("C21").Value = 1

# Spot "C2" (row 31): one more challenger, still not full
$ws.Range("C31").Value = 1

# Spot "E7" (row 50): now full -> Open flips to TRUE, count to 2
$ws.Range("B50").Value = $true
$ws.Range("C50").Value = 2

# Spot "F2" (row 59): one more challenger, still not full
$ws.Range("C59").Value = 1

# Spot "I11" (row 96): now full -> Open flips to TRUE, count to 2
$ws.Range("B96").Value = $true
$ws.Range("C96").Value = 2

# Spot "K9" (row 123): one more challenger, still not full
$ws.Range("C123").Value = 1

# Spot "K10" (row 124): now full -> Open flips to TRUE, count to 2
$ws.Range("B124").Value = $true
$ws.Range("C124").Value = 2

# Spot "M3" (row 145): now full -> Open flips to TRUE, count to 2
$ws.Range("B145").Value = $true
$ws.Range("C145").Value = 2

# Reflect the scrolled viewport / new selection left behind after editing
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
$ws.Range("B51").Select()
